# Insert a new row at position 132 (shifting existing rows 132..190 down to 133..191)
# and populate it with a new weekly price record, per the commit:
#   "Fruta / hortaliza, semanal"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 132. This pushes the
# previous rows 132..190 down to 133..191 and grows the used range to R191.
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row 132 with the new record's data.
$ws.Cells.Item(132, 1).Value = 7
$ws.Cells.Item(132, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(132, 3).Value = "Ñuble"
$ws.Cells.Item(132, 4).Value = 44523
$ws.Cells.Item(132, 5).Value = 16
$ws.Cells.Item(132, 6).Value = 100112008
$ws.Cells.Item(132, 7).Value = "Coliflor"
$ws.Cells.Item(132, 8).Value = "Sin especificar"
$ws.Cells.Item(132, 9).Value = "Primera"
$ws.Cells.Item(132, 10).Value = 240
$ws.Cells.Item(132, 11).Value = 700
$ws.Cells.Item(132, 12).Value = 750
$ws.Cells.Item(132, 13).Value = 725
$ws.Cells.Item(132, 14).Value = "`$/unidad"
$ws.Cells.Item(132, 15).Value = "Región del Maule"
$ws.Cells.Item(132, 16).Value = 725
$ws.Cells.Item(132, 17).Value = 1
$ws.Cells.Item(132, 18).Value = "Hortaliza"

Write-Host "Row 132 inserted and populated"
